# Updates cryptos list values (Price / Volume(1h) columns) as refreshed by the scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($address, $value) {
    $cell = $ws.Range($address)
    # Force a Text number format so Excel keeps the assigned value as a literal
    # string (it would otherwise auto-convert parseable numbers like "228.00").
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # Reset to the default style so no stray style index is left on the cell,
    # matching the original (unstyled) cells in this sheet.
    $cell.Style = "Normal"
}

Set-TextCell "D2" "38.392.06"
Set-TextCell "E2" "  +1.64%  "
Set-TextCell "D3" "2.081.81"
Set-TextCell "E3" "  +2.40%  "
Set-TextCell "E4" "  -0.11%  "
Set-TextCell "D5" "228.00"
Set-TextCell "E5" "  +0.04%  "
Set-TextCell "E6" "  +0.69%  "
Set-TextCell "D7" "60.43"
Set-TextCell "E7" "  +0.61%  "
Set-TextCell "E8" "  -0.11%  "
Set-TextCell "D9" "0.382"
Set-TextCell "E9" "  +2.02%  "
Set-TextCell "E10" "  +0.78%  "
Set-TextCell "E11" "  -0.39%  "
Set-TextCell "D12" "2.388.93"
Set-TextCell "E12" "  +2.29%  "
Set-TextCell "E13" "  +2.42%  "
Set-TextCell "D14" "22.40"
Set-TextCell "E14" "  +7.03%  "
Set-TextCell "D15" "0.782"
Set-TextCell "E15" "  +1.71%  "
Set-TextCell "E16" "  +3.49%  "
Set-TextCell "D17" "2.086.49"
Set-TextCell "E17" "  +2.76%  "
Set-TextCell "D18" "38.301.80"
Set-TextCell "E18" "  +1.50%  "
Set-TextCell "D19" "71.71"
Set-TextCell "E19" "  +3.40%  "
Set-TextCell "D20" "6.02"
Set-TextCell "E20" "  +2.06%  "
Set-TextCell "D21" "0.0₃0832"
Set-TextCell "E21" "  +1.54%  "
Set-TextCell "D22" "225.30"
Set-TextCell "E22" "  +0.59%  "
Set-TextCell "E23" "  -0.03%  "
Set-TextCell "E25" "  +2.11%  "
Set-TextCell "D26" "169.89"
Set-TextCell "E26" "  +1.13%  "
Set-TextCell "D27" "9.40"
Set-TextCell "E27" "  +1.03%  "
Set-TextCell "D28" "0.137"
Set-TextCell "E28" "  +6.45%  "
Set-TextCell "D29" "19.03"
Set-TextCell "E29" "  +1.82%  "
Set-TextCell "D30" "1.36"
Set-TextCell "E30" "  +8.25%  "
Set-TextCell "E31" "  -0.30%  "
Set-TextCell "E32" "  +4.97%  "
Set-TextCell "D33" "4.81"
Set-TextCell "E33" "  +7.71%  "
Set-TextCell "D34" "4.50"
Set-TextCell "E34" "  +2.90%  "
Set-TextCell "E35" "  +0.36%  "
Set-TextCell "E36" "  +2.43%  "
Set-TextCell "D37" "6.35"
Set-TextCell "E37" "  -2.60%  "
Set-TextCell "E38" "  +4.82%  "
Set-TextCell "E39" "  +0.03%  "
Set-TextCell "D40" "18.28"
Set-TextCell "E40" "  +2.44%  "
Set-TextCell "D41" "1.539.83"
Set-TextCell "E41" "  +1.04%  "
Set-TextCell "D42" "100.19"
Set-TextCell "E42" "  +3.35%  "
Set-TextCell "D43" "0.0219"
Set-TextCell "E43" "  +2.09%  "
Set-TextCell "D44" "0.0923"
Set-TextCell "E44" "  +1.93%  "
Set-TextCell "E45" "  -1.36%  "
Set-TextCell "E46" "  +8.30%  "
Set-TextCell "D47" "4.11"
Set-TextCell "E47" "  -1.12%  "
Set-TextCell "E48" "  +0.87%  "
Set-TextCell "E49" "  +2.57%  "
Set-TextCell "E50" "  +0.97%  "
Set-TextCell "D51" "2.279.23"
Set-TextCell "E51" "  +2.43%  "
